$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple D-column price updates ---
Set-TextValue $ws.Range('D2') '275.59'
Set-TextValue $ws.Range('D3') '22.92'
Set-TextValue $ws.Range('D4') '6.343'
Set-TextValue $ws.Range('D5') '0.06234'
Set-TextValue $ws.Range('D7') '6.706'
Set-TextValue $ws.Range('D8') '1.368'
Set-TextValue $ws.Range('D9') '0.8383'
Set-TextValue $ws.Range('D11') '0.1640'
Set-TextValue $ws.Range('D12') '0.08353'
Set-TextValue $ws.Range('D13') '0.03363'
Set-TextValue $ws.Range('D14') '0.03104'
Set-TextValue $ws.Range('D40') '0.04697'
Set-TextValue $ws.Range('D41') '0.007027'
Set-TextValue $ws.Range('D42') '0.1169'
Set-TextValue $ws.Range('D43') '0.003349'
Set-TextValue $ws.Range('D44') '0.01260'
Set-TextValue $ws.Range('D45') '0.00006269'
Set-TextValue $ws.Range('D47') '0.8999'
Set-TextValue $ws.Range('D48') '0.03186'
Set-TextValue $ws.Range('D49') '0.00002300'

# --- Row block updates (rows 15-26): ranking reshuffle ---
Set-TextValue $ws.Range('B15') 'ProBitToken'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D15') '0.1270'
Set-TextValue $ws.Range('E15') '14ProBitTokenPROB'

Set-TextValue $ws.Range('B16') 'BitMartToken'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D16') '0.09330'
Set-TextValue $ws.Range('E16') '15BitMartTokenBMX'

Set-TextValue $ws.Range('B17') 'MCDex'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D17') '3.884'
Set-TextValue $ws.Range('E17') '16MCDexMCB'

Set-TextValue $ws.Range('B18') 'BitForexToken'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D18') '0.001648'
Set-TextValue $ws.Range('E18') '17BitForexTokenBF'

Set-TextValue $ws.Range('B19') 'CoinExToken'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D19') '0.04803'
Set-TextValue $ws.Range('E19') '18CoinExTokenCET'

Set-TextValue $ws.Range('B20') 'TigerCash'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D20') '0.006208'
Set-TextValue $ws.Range('E20') '19TigerCashTCH'

Set-TextValue $ws.Range('B21') 'HotbitToken'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D21') '0.005572'
Set-TextValue $ws.Range('E21') '20HotbitTokenHTBWorstin24h'

Set-TextValue $ws.Range('B22') 'BitKan'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D22') '0.001088'
Set-TextValue $ws.Range('E22') '21BitKanKAN'

Set-TextValue $ws.Range('B23') 'NitroEx'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws.Range('D23') '0.0001500'
Set-TextValue $ws.Range('E23') '22NitroExNTX'

Set-TextValue $ws.Range('B24') 'LEO'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D24') '3.730'
Set-TextValue $ws.Range('E24') '23LEOLEO'

Set-TextValue $ws.Range('B25') 'BTSEToken'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D25') '2.359'
Set-TextValue $ws.Range('E25') '24BTSETokenBTSE'

Set-TextValue $ws.Range('B26') 'BitpandaEcosystemToken'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D26') '0.3379'
Set-TextValue $ws.Range('E26') '25BitpandaEcosystemTokenBEST'
